# Update the "想去人数" (desired-attendance) counts in column F across all
# four sheets, per the "gh-pages output generated at 456a3b4" refresh.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (aggregate of the above)

# 展览
$ws1.Range("F3").Value = 27032
$ws1.Range("F5").Value = 647
$ws1.Range("F6").Value = 190
$ws1.Range("F7").Value = 567
$ws1.Range("F9").Value = 375
$ws1.Range("F15").Value = 483
$ws1.Range("F17").Value = 1620
$ws1.Range("F18").Value = 249
$ws1.Range("F19").Value = 601
$ws1.Range("F21").Value = 457
$ws1.Range("F22").Value = 9

# 演出
$ws2.Range("F2").Value = 4524
$ws2.Range("F6").Value = 209
$ws2.Range("F7").Value = 209
$ws2.Range("F11").Value = 455
$ws2.Range("F20").Value = 30

# 本地生活
$ws3.Range("F2").Value = 5178

# 全部类型
$ws4.Range("F2").Value = 47
$ws4.Range("F3").Value = 5178
$ws4.Range("F5").Value = 27032
$ws4.Range("F6").Value = 4524
$ws4.Range("F9").Value = 647
$ws4.Range("F12").Value = 190
$ws4.Range("F13").Value = 209
$ws4.Range("F14").Value = 209
$ws4.Range("F18").Value = 455
$ws4.Range("F19").Value = 567
$ws4.Range("F23").Value = 375
$ws4.Range("F32").Value = 483
$ws4.Range("F35").Value = 1620
$ws4.Range("F36").Value = 249
$ws4.Range("F37").Value = 601
$ws4.Range("F40").Value = 457
$ws4.Range("F41").Value = 9
$ws4.Range("F45").Value = 30
